$d = $word.ActiveDocument

# The last paragraph in the document body is the "Map being a..." list item.
$last = $d.Paragraphs.Last

$r1 = $last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Interoperability with the standard java collections."

$r2 = $p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Different strategies for the Sequence."
